$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.678.02"
$ws.Range("E2").Value = "  +1.17%  "

$ws.Range("D3").Value = "2.210.74"
$ws.Range("E3").Value = "  +0.51%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'291.86"
$ws.Range("E5").Value = "  -1.04%  "

$ws.Range("D6").Value = "'86.25"
$ws.Range("E6").Value = "  +6.40%  "

$ws.Range("E7").Value = "  +1.08%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "'0.470"
$ws.Range("E9").Value = "  +0.58%  "

$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.0784"
$ws.Range("E10").Value = "  +2.12%  "

$ws.Range("B11").Value = "Avalanche"
$ws.Range("C11").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D11").Value = "'30.15"
$ws.Range("E11").Value = "  +3.78%  "

$ws.Range("D12").Value = "'47.39"
$ws.Range("E12").Value = "  +1.35%  "

$ws.Range("D13").Value = "'0.108"
$ws.Range("E13").Value = "  +1.67%  "

$ws.Range("D14").Value = "'6.32"
$ws.Range("E14").Value = "  +1.70%  "

$ws.Range("D15").Value = "2.552.51"
$ws.Range("E15").Value = "  +0.93%  "

$ws.Range("D16").Value = "'14.00"
$ws.Range("E16").Value = "  +0.62%  "

$ws.Range("D17").Value = "2.205.15"
$ws.Range("E17").Value = "  +0.26%  "

$ws.Range("D18").Value = "'0.725"
$ws.Range("E18").Value = "  +2.26%  "

$ws.Range("D19").Value = "39.649.73"
$ws.Range("E19").Value = "  +1.35%  "

$ws.Range("D20").Value = "'11.49"
$ws.Range("E20").Value = "  +12.21%  "

$ws.Range("D21").Value = "0.0₃0878"
$ws.Range("E21").Value = "  +0.99%  "

$ws.Range("D22").Value = "'5.79"
$ws.Range("E22").Value = "  +1.62%  "

$ws.Range("D23").Value = "'65.79"
$ws.Range("E23").Value = "  +1.92%  "

$ws.Range("D24").Value = "'235.33"
$ws.Range("E24").Value = "  +4.20%  "

$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("D26").Value = "'2.46"
$ws.Range("E26").Value = "  +2.65%  "

$ws.Range("D27").Value = "'1.83"
$ws.Range("E27").Value = "  +1.68%  "

$ws.Range("D28").Value = "'22.67"
$ws.Range("E28").Value = "  +0.83%  "

$ws.Range("E29").Value = "  +1.47%  "

$ws.Range("D30").Value = "'9.25"
$ws.Range("E30").Value = "  +2.51%  "

$ws.Range("D31").Value = "'32.66"
$ws.Range("E31").Value = "  +3.62%  "

$ws.Range("D32").Value = "'152.11"
$ws.Range("E32").Value = "  +1.87%  "

$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("D34").Value = "'4.92"
$ws.Range("E34").Value = "  +2.75%  "

$ws.Range("D35").Value = "'0.0715"
$ws.Range("E35").Value = "  +3.14%  "

$ws.Range("E36").Value = "  +1.80%  "

$ws.Range("E37").Value = "  +2.12%  "

$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").Value = "'16.00"
$ws.Range("E38").Value = "  +4.66%  "

$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").Value = "'2.79"
$ws.Range("E39").Value = "  +6.35%  "

$ws.Range("D40").Value = "'0.0984"
$ws.Range("E40").Value = "  +2.69%  "

$ws.Range("D41").Value = "'1.69"
$ws.Range("E41").Value = "  +2.76%  "

$ws.Range("D42").Value = "2.075.84"
$ws.Range("E42").Value = "  +9.28%  "

$ws.Range("D43").Value = "'3.77"
$ws.Range("E43").Value = "  +4.90%  "

$ws.Range("D44").Value = "'2.16"
$ws.Range("E44").Value = "  +5.95%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'10.01"
$ws.Range("E45").Value = "  +11.28%  "

$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0268"
$ws.Range("E46").Value = "  +3.44%  "

$ws.Range("D47").Value = "'17.54"
$ws.Range("E47").Value = "  +9.72%  "

$ws.Range("D48").Value = "'2.60"
$ws.Range("E48").Value = "  +0.05%  "

$ws.Range("D49").Value = "2.425.10"
$ws.Range("E49").Value = "  +0.80%  "

$ws.Range("D50").Value = "'70.70"
$ws.Range("E50").Value = "  -0.77%  "

$ws.Range("D51").Value = "'88.88"
$ws.Range("E51").Value = "  +1.96%  "
